$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.181.47"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").Value = "2.247.64"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'307.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").Value = "'95.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "

$ws.Range("D7").Value = "'0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("D10").Value = "'35.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("D11").Value = "'0.0818"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "'7.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").Value = "2.591.03"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "2.259.67"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").Value = "'0.835"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").Value = "'13.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("D18").Value = "44.070.66"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "0.0₃0977"
$ws.Range("E19").Value = "  +1.70%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.67%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").Value = "'65.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.06%  "

$ws.Range("D23").Value = "'237.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").Value = "'2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").Value = "'9.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("D29").Value = "'37.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.90%  "

$ws.Range("D30").Value = "'6.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").Value = "'20.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("D32").Value = "'152.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.43%  "

$ws.Range("D33").Value = "'0.0806"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.47%  "

$ws.Range("E34").Value = "  +3.55%  "

$ws.Range("E35").Value = "  -2.89%  "

$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("E37").Value = "  +2.31%  "

$ws.Range("D38").Value = "'1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.35%  "

$ws.Range("D39").Value = "'3.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.23%  "

$ws.Range("D40").Value = "'3.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.13%  "

$ws.Range("D41").Value = "'14.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.78%  "

$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "1.742.48"
$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").Value = "'83.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.72%  "

$ws.Range("D46").Value = "'0.191"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.86%  "

$ws.Range("D47").Value = "'100.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("D48").Value = "'4.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.52%  "

$ws.Range("D49").Value = "'8.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").Value = "'54.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "

$ws.Range("D51").Value = "'68.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.06%  "
